$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project Planner")

# Clear the now-unused column header labels (PLAN START, PLAN DURATION,
# ACTUAL START, ACTUAL DURATION, PERCENT COMPLETE) while keeping their
# formatting intact.
$ws.Range("C3:G3").Value = ""

# Leave the selection on the merged cell E3:E4 (the last-edited header cell).
$ws.Range("E3:E4").Select()
